# Append a new job-listing row to the "ランサーズ" sheet and refresh the
# "取得日時" (fetched-at) timestamp for every existing row, per the
# 2026-01-30 06:56 JST scrape run.
#
# Net effect vs. the previous state:
#   - A new row is inserted at row 3 (pushing the former rows 3-6 down to 4-7)
#   - Every row's column A timestamp becomes "2026-01-30 06:56:13"
#   - The hyperlinks in column F are rebuilt so each one points at the URL
#     that now actually lives in that row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-30 06:56:13"

# Drop the existing hyperlink objects now (their cell anchors would
# otherwise stay pinned to old cell addresses once the new row is inserted).
$ws.Hyperlinks.Delete()

# Make room for the new listing right under row 2 (the highest scoring,
# unchanged, Java job). This shifts the old rows 3-6 down to rows 4-7.
$ws.Rows.Item(3).Insert()

# --- Refresh the "fetched at" timestamp on every data row ---
$ws.Range("A2").Value = $newTimestamp
$ws.Range("A3").Value = $newTimestamp
$ws.Range("A4").Value = $newTimestamp
$ws.Range("A5").Value = $newTimestamp
$ws.Range("A6").Value = $newTimestamp
$ws.Range("A7").Value = $newTimestamp

# --- Fill in the brand-new listing in row 3 ---
$ws.Range("B3").Value = "店舗タブレット用Webサイトと管理画面の制作依頼"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5482389"
$ws.Range("G3").Value = 63
$ws.Range("H3").Value = "◇サイト"

# --- The rows that used to be 4 and 5 gain/keep their skill-summary tag ---
$ws.Range("H5").Value = "◇管理"

# --- Re-create the hyperlinks so column F again matches its displayed URL ---
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5482097")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5482389")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5481859")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5418064")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5481715")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5481888")

# Adding hyperlinks re-applies a fresh "Hyperlink" cell style; re-stamp the
# named style so these cells reuse the workbook's existing Hyperlink xf
# instead of accumulating an unused duplicate.
$ws.Range("F2:F7").Style = "Hyperlink"
